$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 377, pushing existing rows 377-392 down to 378-393
$ws.Rows("377:377").Insert()

# Populate the newly inserted row 377 with the new weekly price record
$ws.Range("A377").Value = 1
$ws.Range("B377").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C377").Value = "Arica y Parinacota"
$ws.Range("D377").Value = 45041
$ws.Range("E377").Value = 15
$ws.Range("F377").Value = "Fruta"
$ws.Range("G377").Value = 100108
$ws.Range("H377").Value = "Tropicales y subtropicales"
$ws.Range("I377").Value = 100108006
$ws.Range("J377").Value = "Plátano"
$ws.Range("K377").Value = "Sin especificar"
$ws.Range("L377").Value = "Pintón"
$ws.Range("M377").Value = 110
$ws.Range("N377").Value = 17000
$ws.Range("O377").Value = 18000
$ws.Range("P377").Value = 17455
$ws.Range("Q377").Value = "`$/caja 20 kilos"
$ws.Range("R377").Value = "Ecuador"
$ws.Range("S377").Value = 873
$ws.Range("T377").Value = 20
